{"js": "// Replace the date line and each two-digit multiplication problem's\n// text with its new value. Every \"find\" string is unique in the\n// document, so a matchCase/matchWholeWord search + Replace is safe.\nconst replacements = [\n  [\"2025-06-21 Saturday\", \"2025-06-22 Sunday\"],\n  [\"67\u00d788=\", \"18\u00d761=\"],\n  [\"41\u00d784=\", \"37\u00d753=\"],\n  [\"40\u00d766=\", \"34\u00d718=\"],\n  [\"59\u00d760=\", \"46\u00d712=\"],\n  [\"49\u00d765=\", \"91\u00d758=\"],\n  [\"29\u00d798=\", \"80\u00d787=\"],\n  [\"89\u00d752=\", \"81\u00d750=\"],\n  [\"99\u00d757=\", \"61\u00d744=\"],\n  [\"87\u00d748=\", \"52\u00d741=\"],\n  [\"29\u00d774=\", \"49\u00d758=\"],\n  [\"22\u00d741=\", \"28\u00d772=\"],\n  [\"96\u00d726=\", \"73\u00d749=\"],\n  [\"74\u00d750=\", \"60\u00d740=\"],\n  [\"33\u00d737=\", \"43\u00d765=\"],\n  [\"72\u00d763=\", \"36\u00d728=\"],\n  [\"71\u00d741=\", \"53\u00d719=\"],\n  [\"51\u00d792=\", \"29\u00d780=\"],\n  [\"38\u00d718=\", \"67\u00d754=\"],\n  [\"46\u00d798=\", \"35\u00d771=\"],\n  [\"87\u00d717=\", \"38\u00d778=\"],\n  [\"76\u00d716=\", \"53\u00d780=\"],\n  [\"49\u00d775=\", \"93\u00d727=\"],\n  [\"42\u00d738=\", \"16\u00d791=\"],\n  [\"85\u00d753=\", \"88\u00d725=\"],\n  [\"99\u00d774=\", \"96\u00d795=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [findText, newText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${findText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each two-digit multiplication problem's\n# text with its new value. Every \"find\" string is unique in the\n# document, so a simple Find/Replace (wdReplaceAll, matching exactly\n# one occurrence each) is safe and idempotent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-21 Saturday\", \"2025-06-22 Sunday\"),\n    @(\"67\u00d788=\", \"18\u00d761=\"),\n    @(\"41\u00d784=\", \"37\u00d753=\"),\n    @(\"40\u00d766=\", \"34\u00d718=\"),\n    @(\"59\u00d760=\", \"46\u00d712=\"),\n    @(\"49\u00d765=\", \"91\u00d758=\"),\n    @(\"29\u00d798=\", \"80\u00d787=\"),\n    @(\"89\u00d752=\", \"81\u00d750=\"),\n    @(\"99\u00d757=\", \"61\u00d744=\"),\n    @(\"87\u00d748=\", \"52\u00d741=\"),\n    @(\"29\u00d774=\", \"49\u00d758=\"),\n    @(\"22\u00d741=\", \"28\u00d772=\"),\n    @(\"96\u00d726=\", \"73\u00d749=\"),\n    @(\"74\u00d750=\", \"60\u00d740=\"),\n    @(\"33\u00d737=\", \"43\u00d765=\"),\n    @(\"72\u00d763=\", \"36\u00d728=\"),\n    @(\"71\u00d741=\", \"53\u00d719=\"),\n    @(\"51\u00d792=\", \"29\u00d780=\"),\n    @(\"38\u00d718=\", \"67\u00d754=\"),\n    @(\"46\u00d798=\", \"35\u00d771=\"),\n    @(\"87\u00d717=\", \"38\u00d778=\"),\n    @(\"76\u00d716=\", \"53\u00d780=\"),\n    @(\"49\u00d775=\", \"93\u00d727=\"),\n    @(\"42\u00d738=\", \"16\u00d791=\"),\n    @(\"85\u00d753=\", \"88\u00d725=\"),\n    @(\"99\u00d774=\", \"96\u00d795=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceAll = 2 (named Word constants aren't\n    # pre-seeded as PowerShell variables in this host, so use the literal\n    # values Word's object model defines for them).\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n$d.Save()\n"}
